# Weekly update: insert a new record (row) for "Vega Modelo de Temuco - Espárragos"
# This pushes the existing rows 61-71 down to 62-72 and adds a brand new
# observation as the new row 61 (dated 2022-10-17 / serial 44841).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 61 (rows 61..71 shift down to 62..72)
$ws.Rows.Item(61).Insert()

# Populate the new row 61 with the new observation
$ws.Cells.Item(61, 1).Value  = 10
$ws.Cells.Item(61, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(61, 3).Value  = 'La Araucanía'
$ws.Cells.Item(61, 4).Value  = 44841
$ws.Cells.Item(61, 5).Value  = 9
$ws.Cells.Item(61, 6).Value  = 300000000
$ws.Cells.Item(61, 7).Value  = 'Espárragos'
$ws.Cells.Item(61, 8).Value  = 'Sin especificar'
$ws.Cells.Item(61, 9).Value  = 'Primera'
$ws.Cells.Item(61, 10).Value = 200
$ws.Cells.Item(61, 11).Value = 1700
$ws.Cells.Item(61, 12).Value = 1700
$ws.Cells.Item(61, 13).Value = 1700
$ws.Cells.Item(61, 14).Value = '$/kilo'
$ws.Cells.Item(61, 15).Value = 'Región del Maule'
$ws.Cells.Item(61, 16).Value = 1700
$ws.Cells.Item(61, 17).Value = 1
$ws.Cells.Item(61, 18).Value = 'Hortaliza'
